# Refresh the leve-profit market-data columns (H:N) for the rows
# touched by this run, one worksheet-row at a time. Values come from
# the scheduled market-price refresh; cells that no longer carry a
# computed figure for a row are cleared outright.
$wb = $excel.ActiveWorkbook

# ALC!100 - Asking for a Friend / Beetle Glue
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 967.3333
$ws.Range("I100").Value = 983.2
$ws.Range("J100").Value = 888
$ws.Range("K100").Value = 983.2
$ws.Range("L100").Value = 888
$ws.Range("M100").Value = -442.2
$ws.Range("N100").Value = -1970

# ARM!32 - Ingot We Trust / Steel Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10842.096
$ws.Range("I32").Value = 5223.769
$ws.Range("K32").Value = 5223.769
$ws.Range("M32").Value = -4936.769

# ARM!45 - Hollow Hallmarks / Mythril Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3126.875
$ws.Range("I45").Value = 2516.625
$ws.Range("K45").Value = 2516.625
$ws.Range("M45").Value = -2139.625

# ARM!61 - Dealing with the Tough Stuff / Cobalt Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2823.76
$ws.Range("I61").Value = 2022.8889
$ws.Range("J61").Value = 4883.143
$ws.Range("K61").Value = 2022.8889
$ws.Range("L61").Value = 4883.143
$ws.Range("M61").Value = -1810.8889
$ws.Range("N61").Value = -5307.143

# ARM!74 - As the Bolt Flies / Titanium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1391.9032
$ws.Range("I74").Value = 805.1111
$ws.Range("K74").Value = 805.1111
$ws.Range("M74").Value = 68.88890000000004

# ARM!77 - Heavy Metal Banned (L) / Titanium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1391.9032
$ws.Range("I77").Value = 805.1111
$ws.Range("K77").Value = 4025.5555
$ws.Range("M77").Value = 342.4445000000001

# ARM!136 - Metal with Mettle / Cobalt Tungsten Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2823.76
$ws.Range("I136").Value = 2022.8889
$ws.Range("J136").Value = 4883.143
$ws.Range("K136").Value = 6068.6667
$ws.Range("L136").Value = 14649.429
$ws.Range("M136").Value = -3518.6667
$ws.Range("N136").Value = -19749.429

# BSM!134 - Ruthenium Supremium / Ruthenium Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3473.2104
$ws.Range("I134").Value = 3078.6
$ws.Range("K134").Value = 9235.799999999999
$ws.Range("M134").Value = -6700.799999999999

# CRP!122 - Timber of Tenkonto / Horse Chestnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 5000
$ws.Range("K122").Value = 15000
$ws.Range("M122").Value = -12550

# CUL!7 - It's Always Sunny in Vylbrand / Raisins
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 237.66667
$ws.Range("I7").Value = 275.2
$ws.Range("K7").Value = 825.5999999999999
$ws.Range("M7").Value = -713.5999999999999

# CUL!92 - Oh No Udon / Gyr Abanian Flour
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 101
$ws.Range("J92").Value = 101
$ws.Range("L92").Value = 303
$ws.Range("N92").Value = -2799

# CUL!107 - Slippery Service / Frantoio Oil
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 740.34784
$ws.Range("I107").Value = 645.75
$ws.Range("J107").Value = 760.2632
$ws.Range("K107").Value = 1937.25
$ws.Range("L107").Value = 2280.7896
$ws.Range("M107").Value = -17.25
$ws.Range("N107").Value = -6120.7896

# CUL!137 - Creative Chocolate / Gateau au Chocolat
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 1503.6364
$ws.Range("I137").Value = 1503.6364
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 4510.9092
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 589.0907999999999
$ws.Range("N137").ClearContents()

# GSM!70 - Sky Is the Limit / Mythrite Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6120.7144
$ws.Range("J70").Value = 7565.5557
$ws.Range("L70").Value = 7565.5557
$ws.Range("N70").Value = -8105.5557

# GSM!73 - Hulls of Broken Dreams (L) / Mythrite Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6120.7144
$ws.Range("J73").Value = 7565.5557
$ws.Range("L73").Value = 7565.5557
$ws.Range("N73").Value = -9437.555700000001

# GSM!123 - Workplace Workout / Ametrine Ring of Fending
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 31432.27
$ws.Range("I123").Value = 28190.578
$ws.Range("J123").Value = 40231.145
$ws.Range("K123").Value = 28190.578
$ws.Range("L123").Value = 40231.145
$ws.Range("M123").Value = -25740.578
$ws.Range("N123").Value = -45131.145

# GSM!132 - On Board for Lar / Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3140.5454
$ws.Range("I132").Value = 2328.3845
$ws.Range("K132").Value = 6985.1535
$ws.Range("M132").Value = -4455.1535

# LTW!7 - Tan Before the Ban / Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

# LTW!22 - Skin off Their Backs / Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1953.8788
$ws.Range("I22").Value = 867.63635
$ws.Range("K22").Value = 867.63635
$ws.Range("M22").Value = -572.63635

# LTW!27 - Fire and Hide / Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1953.8788
$ws.Range("I27").Value = 867.63635
$ws.Range("K27").Value = 867.63635
$ws.Range("M27").Value = -760.63635

# LTW!40 - Best Served Toad / Toad Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4432.091
$ws.Range("I40").Value = 3425
$ws.Range("J40").Value = 6194.5
$ws.Range("K40").Value = 3425
$ws.Range("L40").Value = 6194.5
$ws.Range("M40").Value = -3289
$ws.Range("N40").Value = -6466.5

# LTW!55 - It's Not a Job, It's a Calling / Peiste Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 282.57895
$ws.Range("I55").Value = 302.16666
$ws.Range("J55").Value = 249
$ws.Range("K55").Value = 302.16666
$ws.Range("L55").Value = 249
$ws.Range("M55").Value = -129.16666
$ws.Range("N55").Value = -595

# LTW!61 - Spelling Me Softly / Raptor Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3735.6667
$ws.Range("I61").Value = 3735.6667
$ws.Range("K61").Value = 3735.6667
$ws.Range("M61").Value = -3533.6667

# LTW!93 - Hide to Go Seek / Gagana Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3862.5
$ws.Range("I93").Value = 3862.5
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 3862.5
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -2614.5
$ws.Range("N93").ClearContents()

# LTW!100 - Tiger in the Sack / Tiger Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2871.2307
$ws.Range("I100").Value = 2755.6667
$ws.Range("K100").Value = 2755.6667
$ws.Range("M100").Value = -2214.6667

# LTW!113 - Peace in Rest / Atrociraptor Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3735.6667
$ws.Range("I113").Value = 3735.6667
$ws.Range("K113").Value = 3735.6667
$ws.Range("M113").Value = -1565.6667

# LTW!122 - Hell on Leather / Gaja Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4928.3076
$ws.Range("I122").Value = 4274.5557
$ws.Range("K122").Value = 12823.6671
$ws.Range("M122").Value = -10373.6671

# LTW!126 - Battered Books / Saiga Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# WVR!122 - Heavy Armoire / Dark Hempen Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3205.8948
$ws.Range("I122").Value = 2945.25
$ws.Range("J122").Value = 3395.4546
$ws.Range("K122").Value = 8835.75
$ws.Range("L122").Value = 10186.3638
$ws.Range("M122").Value = -6385.75
$ws.Range("N122").Value = -15086.3638

# WVR!126 - A Polished Purchase / Snow Linen
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4500
$ws.Range("I126").Value = 4500
$ws.Range("K126").Value = 13500
$ws.Range("M126").Value = -11030

# WVR!132 - Comfy Cabins / Snow Cotton Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5081
$ws.Range("I132").Value = 5575.647
$ws.Range("K132").Value = 16726.941
$ws.Range("M132").Value = -14196.941
